$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 85
$wsExhibition.Range("F4").Value = 2233
$wsExhibition.Range("F6").Value = 370

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 85
$wsAll.Range("F4").Value = 2233
$wsAll.Range("F7").Value = 370
